$wb = $excel.ActiveWorkbook

$wsParameter = $wb.Worksheets.Item("Parameter")
$wsFilepath  = $wb.Worksheets.Item("Filepath")

# Update the custom boolean number format used by Parameter!B7:B12
$wsParameter.Range("B7:B12").NumberFormat = '"BOOL"e"AN"'

# Remove "saving generation" and "saving exchanges" parameters from Filepath sheet.
# Row 7 ("saving generation") contents are cleared, row 8 ("saving exchanges")
# and the blank spacer row 9 are fully deleted, shifting subsequent rows up.
$wsFilepath.Range("A7:B7").ClearContents()
$wsFilepath.Rows("8:9").Delete()

# Make Parameter the active sheet/tab (was Filepath).
$wsParameter.Select()

$wb.Save()
